# Update the clinical-variables list in column D:
#  - insert "Sample_pretreated" right after "Pembro_status"
#  - rename "PDL1_percent_score_group" to "PDL1_IHC_percent_DrYang_group"
#  - include the previously-missing "OS_months_group_quartile" entry
# This shifts every entry from the old D3 onward down by one row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Write the PDL1 rename first (new shared-string entry), then the rest of the
# list including the newly-inserted "Sample_pretreated" row, so new strings
# land in the workbook's shared-string table in the same order the author's
# edit produced them.
$ws.Range("D22").Value = "PDL1_IHC_percent_DrYang_group"

$clinicalVars = @(
    "Pembro_status",
    "Sample_pretreated",
    "Sex",
    "Age_at_diagnosis_abbrev",
    "Type",
    "Sample_type",
    "Body_part",
    "Body_part_abbrev",
    "Metastasis_brain",
    "Metastasis_brain_status",
    "Smoking_status_abbrev",
    "Patient_history_of_cancer_abbrev",
    "Family_history_of_cancer_abbrev",
    "predictive_biomarker",
    "Best_overall_response_group",
    "Best_overall_response_detailed",
    "PFS_months_group_median",
    "PFS_months_group_quartile",
    "OS_months_group_median",
    "OS_months_group_quartile",
    "PDL1_IHC_percent_DrYang_group",
    "Impact_TMB_score_group10",
    "STK11",
    "KEAP1"
)

$row = 2
foreach ($name in $clinicalVars) {
    $ws.Range("D$row").Value = $name
    $row = $row + 1
}

# Column D now holds a longer string ("PDL1_IHC_percent_DrYang_group"), so widen
# it to fit, matching the workbook's existing best-fit column sizing.
$ws.Columns.Item(4).ColumnWidth = 28.67
